$d = $word.ActiveDocument
Write-Output "content end: $($d.Content.End)"
foreach ($pos in 4500,4540,4545,4547,4550,4600,4605,4606,4607) {
  try {
    $r2 = $d.Range($pos, $pos)
    Write-Output "pos=$pos range=[$($r2.Start),$($r2.End)]"
  } catch {
    Write-Output "pos=$pos ERR $_"
  }
}
